# Update citation placeholder references in each paragraph.
# Note: "Ref-DJ49F2" occurs twice in the original text (paragraphs 2 and 5)
# but must map to two different new values, so each replacement is scoped
# to its own paragraph's Range rather than run globally over $d.Content.

$d = $word.ActiveDocument

$replacements = @(
    @{ Index = 1; Old = "Ref-AB1CD2"; New = "Ref-f819506" },
    @{ Index = 2; Old = "Ref-DJ49F2"; New = "Ref-s842184" },
    @{ Index = 3; Old = "Ref-J7X8N2"; New = "Ref-u470921" },
    @{ Index = 4; Old = "Ref-A1B2C3"; New = "Ref-f217253" },
    @{ Index = 5; Old = "Ref-DJ49F2"; New = "Ref-f652463" },
    @{ Index = 6; Old = "Ref-AB12CD"; New = "Ref-s229328" }
)

foreach ($r in $replacements) {
    $rng = $d.Paragraphs($r.Index).Range
    $rng.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
